$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$strings = @(
    "自動実行",
    "掃除道具入れ",
    "石鹸だ・・・`n残り少なくなっている・・・",
    "EV003",
    "--------キャラ指定-----座標設定------------",
    "----------------------------------------",
    "何か使えるものは・・・",
    "食糧庫とタブの付いた鍵を見つけた！`n他には・・・・",
    "シィナのパンツを見つけてしまった・・・",
    "\C[10]欲しい",
    "嫌な予感がする",
    "・・・`nごくり・・・",
    "これを拾うと色々と終わってしまう気がする・・・`n戻しておこう。",
    "\n<\n[3]>それアタシのパンツにゃん。`nそれ拾ってどーするつもりにゃ？",
    "\n<\n[3]>ふーん・・・`nどーせくんくんしながらオナニーするにゃ？`nどすけべマゾ。",
    "\n<\n[3]>いいにゃ。`nそれあげるにゃん。`nプレゼントにゃ♥",
    "\n<\n[3]>ただしお前じゃなくてお前のチンポにな！",
    "\n<\n[3]>ただしお前じゃなくてお前のチンポににゃ！",
    "炎",
    " <enemy:99><ch:1.5>",
    "灯り 5 3 255 0.3",
    "\n<\n[1]>あっつ！！！",
    "変数203（ARGP攻撃種類）`n1斬　2打撃　3水　4火　5雷`n6誘惑　7食べ物`n特殊206（個別攻撃種類）`n1リンゴ　2皿　3卵",
    " <enemy:99><cw:1.5><ch:1.5>",
    "立体起動",
    "<enemy:99>",
    "水たまりサンプル",
    "食糧庫ドア",
    "掃除道具入れだ・・・`n何か使えるものがあるかもしれない・・・",
    "調べる",
    "今はやめておく",
    "もう使えそうなものはなさそうだ・・・",
    "モップの棒が取れてしまった・・・`n何かに使えるかもしれない。",
    "\n<リリー>お風呂お風呂～♪",
    "\n<シィナ>アタシも入るにゃ！",
    "\n<ライム>えー、じゃあ私もー♥",
    "\n<リリー>ひとりくらいあいつ探しなさいよ・・・",
    "お風呂の鏡は曇っている・・・",
    "EV014",
    "水が入った桶・・・",
    "湿気でふにゃふにゃになってしまっている・・・",
    "何故お風呂に焼き魚が・・・",
    "何故お風呂にパイが・・・",
    "薪だ。`nこんなところに置いて湿気ないのか。",
    "いい匂いのタオルだ。`nでもあまり嗅がない方がいい・・・",
    "報告書",
    "鏡に自分が映っている・・・",
    "\n<\n[1]>（へんなもの付けられちゃったし・・・`nもう女湯に入れないのかな・・・）",
    "\n<\n[1]>（大浴場・・・`n随分と立派だ。`n余程の金持ちが住んでいたに違いない。）",
    "大きいタライだ。`nもし降ってきたら首を痛めそうなサイズ・・・",
    "シィナ決戦",
    "\n<ライム>\n[1]・・・`n本当に戦わなきゃダメなの？",
    "\n<ライム>・・・`nそっか。`nじゃあしょうがないね。",
    "\n<ライム>どうなっても恨みっこ無しだよ！`n絶対に負けないんだから！",
    "ポップアップ有効化",
    "\n<ライム>はぁ・・・はぁ・・・♥`nと、とけるぅ・・・♥`nあへぇ～・・・♥",
    "ライムを退治した！",
    "リリーの日記",
    "ライムの日記`nお風呂でピザ",
    "読んでみる",
    "やめておく",
    "ライムの日記`n水かけ論",
    "EV039",
    "紙が落ちるような音がした・・・",
    "もう何もない・・・",
    "EV041",
    "災害対策の本だ・・・",
    "\n[1]は\C[3]『警戒』\C[0]を覚えた！",
    "\n<\n[1]>（少しやつれたように見える・・・`nでも、ダイエットと思えば。）",
    "\n<\n[1]>（少しやつれたかもしれない・・・）",
    "EV043",
    "EV044",
    "EV045",
    "EV046",
    "EV047",
    "EV048",
    "EV049",
    "EV050",
    "EV051",
    "EV052",
    "EV053",
    "EV054",
    "EV055",
    "ライムの日記`n脱走者"
)

# Clear any pre-existing content in columns C and D (no longer used)
$ws.Range("C1:D84").ClearContents()

for ($i = 0; $i -lt $strings.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $strings[$i]
}

$onlyARows = @(5, 6, 21, 23)
for ($i = 0; $i -lt $strings.Length; $i++) {
    $row = $i + 1
    if ($onlyARows -notcontains $row) {
        $ws.Cells.Item($row, 2).Value = $strings[$i]
    }
}
